$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.159.80"
$ws.Range("D3").Value = "2.594.14"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.99"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.32"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "2.614.64"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "3.054.33"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "58.136.23"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.37"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "2.566.91"
$ws.Range("E18").Value = "  -3.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.35"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("E22").Value = "  +2.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.39"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.405"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.716.74"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("E30").Value = "  -5.16%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.13"
$ws.Range("E32").Value = "  -5.93%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.78"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.69"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("E37").Value = "  -3.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.866"
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("E40").Value = "  +2.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.02"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "270.39"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0957"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.73"
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("E49").Value = "  -1.43%  "
$ws.Range("D50").Value = "1.967.79"
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("E51").Value = "  +3.07%  "
